$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 469 (shifts existing rows 469-501 down to 470-502),
# copying formatting from the row above (matches Excel's default Insert behaviour,
# which is what produces the dimension A1:R502 and the s="2" date style on column D).
$ws.Rows(469).Insert()

# Populate the newly inserted row 469 with the new weekly price record.
$ws.Cells.Item(469, 1).Value = 7
$ws.Cells.Item(469, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(469, 3).Value = "Ñuble"
$ws.Cells.Item(469, 4).Value = 45166
$ws.Cells.Item(469, 5).Value = 16
$ws.Cells.Item(469, 6).Value = 100112006
$ws.Cells.Item(469, 7).Value = "Repollo"
$ws.Cells.Item(469, 8).Value = "Crespo record"
$ws.Cells.Item(469, 9).Value = "Primera"
$ws.Cells.Item(469, 10).Value = 300
$ws.Cells.Item(469, 11).Value = 1000
$ws.Cells.Item(469, 12).Value = 1200
$ws.Cells.Item(469, 13).Value = 1100
$ws.Cells.Item(469, 14).Value = "$/unidad"
$ws.Cells.Item(469, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(469, 16).Value = 1100
$ws.Cells.Item(469, 17).Value = 1
$ws.Cells.Item(469, 18).Value = "Hortaliza"
